$d = $word.ActiveDocument
$aposR = [char]8217

# The document currently ends with a pre-existing blank paragraph right after
# the "...troubles focusing." entry. We reuse that blank paragraph as the
# blank separator line called for by the new diary entry, then append the
# new "21/06/2025:" heading paragraph, the new diary-entry paragraph, and
# finally a fresh blank paragraph to replace the one we just reused (so the
# document still ends with a blank line, as before).

# --- bold "21/06/2025:" heading paragraph ---------------------------------
$lastRng = $d.Paragraphs($d.Paragraphs.Count).Range
$lastRng.InsertParagraphAfter() | Out-Null

$dateText = "21/06/2025:"
$dateRng = $d.Paragraphs($d.Paragraphs.Count).Range
$dateRng.InsertAfter($dateText) | Out-Null

$dateParaRng = $d.Paragraphs($d.Paragraphs.Count).Range
$dateParaRng.Font.Name = "Times New Roman"
$dateParaRng.Font.NameAscii = "Times New Roman"
$dateParaRng.Font.NameFarEast = "Times New Roman"
$dateParaRng.Font.Bold = 1
$dateParaRng.LanguageID = 1055

# --- diary-entry paragraph --------------------------------------------------
$dateParaRng2 = $d.Paragraphs($d.Paragraphs.Count).Range
$dateParaRng2.InsertParagraphAfter() | Out-Null

$bodyText = "This has been a real unproductive day for me. I just solved a bunch of linear algebra problems from sections 1.1 and 1.2, and I didn" + $aposR + "t write one line of code. I played a bunch of games though, I guess it was just a calm down day for me."
$bodyRng = $d.Paragraphs($d.Paragraphs.Count).Range
$bodyRng.InsertAfter($bodyText) | Out-Null

$bodyParaRng = $d.Paragraphs($d.Paragraphs.Count).Range
$bodyParaRng.Font.Name = "Times New Roman"
$bodyParaRng.Font.NameAscii = "Times New Roman"
$bodyParaRng.Font.NameFarEast = "Times New Roman"
$bodyParaRng.Font.Bold = 0
$bodyParaRng.LanguageID = 1055

# --- trailing blank paragraph (mirrors the document's original ending) ----
$bodyParaRng2 = $d.Paragraphs($d.Paragraphs.Count).Range
$bodyParaRng2.InsertParagraphAfter() | Out-Null

$trailingRng = $d.Paragraphs($d.Paragraphs.Count).Range
$trailingRng.Font.Name = "Times New Roman"
$trailingRng.Font.NameAscii = "Times New Roman"
$trailingRng.Font.NameFarEast = "Times New Roman"
$trailingRng.Font.Bold = 0
$trailingRng.LanguageID = 1055
